$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 21, pushing existing rows 21-24 down to 23-26.
$ws.Rows("21:22").Insert()

# New row 21
$ws.Range("A21").Value = 5
$ws.Range("B21").Value = "Macroferia Regional de Talca"
$ws.Range("C21").Value = "Maule"
$ws.Range("D21").Value = "12/13/2021"
$ws.Range("E21").Value = 7
$ws.Range("F21").Value = "Fruta"
$ws.Range("G21").Value = 100103
$ws.Range("H21").Value = "Frutos de hueso (carozo)"
$ws.Range("I21").Value = 100103003
$ws.Range("J21").Value = "Damasco"
$ws.Range("K21").Value = "Castle Brite"
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 100
$ws.Range("N21").Value = 18000
$ws.Range("O21").Value = 18000
$ws.Range("P21").Value = 18000
$ws.Range("Q21").Value = "$/caja 15 kilos"
$ws.Range("R21").Value = "Región de O'Higgins"
$ws.Range("S21").Value = 1200
$ws.Range("T21").Value = 15

# New row 22
$ws.Range("A22").Value = 5
$ws.Range("B22").Value = "Macroferia Regional de Talca"
$ws.Range("C22").Value = "Maule"
$ws.Range("D22").Value = "12/13/2021"
$ws.Range("E22").Value = 7
$ws.Range("F22").Value = "Fruta"
$ws.Range("G22").Value = 100103
$ws.Range("H22").Value = "Frutos de hueso (carozo)"
$ws.Range("I22").Value = 100103003
$ws.Range("J22").Value = "Damasco"
$ws.Range("K22").Value = "Castle Brite"
$ws.Range("L22").Value = "Segunda"
$ws.Range("M22").Value = 50
$ws.Range("N22").Value = 15000
$ws.Range("O22").Value = 15000
$ws.Range("P22").Value = 15000
$ws.Range("Q22").Value = "$/caja 16 kilos"
$ws.Range("R22").Value = "Región de O'Higgins"
$ws.Range("S22").Value = 1000
$ws.Range("T22").Value = 15
